$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44333
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 58
$ws.Range("D3").Value = 44333
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 9000
$ws.Range("S3").Value = 900
$ws.Range("D4").Value = 44333
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("S4").Value = 800
$ws.Range("D5").Value = 44343
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 47
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("D6").Value = 44343
$ws.Range("L6").Value = 'Primera'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("D7").Value = 44343
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 58
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 800
$ws.Range("D8").Value = 44307
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 40
$ws.Range("D9").Value = 44301
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("S9").Value = 1000
$ws.Range("D10").Value = 44308
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 1000
$ws.Range("D11").Value = 44308
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 48
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 800
$ws.Range("D12").Value = 44326
$ws.Range("M12").Value = 65
$ws.Range("D13").Value = 44326
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 67
$ws.Range("N13").Value = 8000
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 8000
$ws.Range("S13").Value = 800
$ws.Range("D14").Value = 44319
$ws.Range("M14").Value = 68
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 57
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8000
$ws.Range("P15").Value = 8000
$ws.Range("S15").Value = 800
$ws.Range("D16").Value = 44329
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 56
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 9000
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 900
$ws.Range("D17").Value = 44329
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 800
$ws.Range("D18").Value = 44309
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 45
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("S18").Value = 1000
$ws.Range("D19").Value = 44322
$ws.Range("M19").Value = 56
$ws.Range("D20").Value = 44322
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = 8000
$ws.Range("O20").Value = 8000
$ws.Range("P20").Value = 8000
$ws.Range("R20").Value = 'Provincia de Quillota'
$ws.Range("S20").Value = 800
$ws.Range("D21").Value = 44328
$ws.Range("M21").Value = 45
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("R21").Value = 'Provincia de Quillota'
$ws.Range("S21").Value = 800
$ws.Range("D22").Value = 44328
$ws.Range("M22").Value = 48
$ws.Range("N22").Value = 7000
$ws.Range("O22").Value = 7000
$ws.Range("P22").Value = 7000
$ws.Range("R22").Value = 'Provincia de Quillota'
$ws.Range("S22").Value = 700
$ws.Range("D23").Value = 44315
$ws.Range("M23").Value = 45
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 10000
$ws.Range("P23").Value = 10000
$ws.Range("R23").Value = 'Provincia de Quillota'
$ws.Range("S23").Value = 1000
$ws.Range("D24").Value = 44312
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 48
$ws.Range("N24").Value = 10000
$ws.Range("O24").Value = 10000
$ws.Range("P24").Value = 10000
$ws.Range("R24").Value = 'Provincia de Quillota'
$ws.Range("S24").Value = 1000
$ws.Range("D25").Value = 44314
$ws.Range("M25").Value = 47
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 9000
$ws.Range("S25").Value = 900
$ws.Range("D26").Value = 44323
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 60
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("S26").Value = 1000
$ws.Range("D27").Value = 44323
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 9000
$ws.Range("S27").Value = 900
$ws.Range("D28").Value = 44306
$ws.Range("M28").Value = 45
$ws.Range("D29").Value = 44321
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 58
$ws.Range("N29").Value = 9000
$ws.Range("O29").Value = 9000
$ws.Range("P29").Value = 9000
$ws.Range("S29").Value = 900
$ws.Range("D30").Value = 44302
$ws.Range("M30").Value = 45
